# Fruta / hortaliza, semanal
# Adds this week's price-report rows for "Terminal La Palmera de La Serena - Pera"
# (Packham's Triumph, 3 calidades) right above the existing data block, pushing the
# older rows down by 3 (the oldest 3 rows that fall past the end of the table are
# re-appended at the bottom automatically, since the used range grows by 3 rows).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 fresh rows right before the existing row 699, shifting the old
# rows 699:723 down to 702:726 (dimension grows from T723 to T726).
$ws.Rows("699:701").Insert()

function Set-Row {
    param($Row, $A, $B, $C, $D, $E, $F, $G, $H, $I, $J, $K, $L, $M, $N, $O, $P, $Q, $R, $S, $T)

    $ws.Cells.Item($Row, 1).Value = $A
    $ws.Cells.Item($Row, 2).Value = $B
    $ws.Cells.Item($Row, 3).Value = $C
    $ws.Cells.Item($Row, 4).Value = $D
    $ws.Cells.Item($Row, 5).Value = $E
    $ws.Cells.Item($Row, 6).Value = $F
    $ws.Cells.Item($Row, 7).Value = $G
    $ws.Cells.Item($Row, 8).Value = $H
    $ws.Cells.Item($Row, 9).Value = $I
    $ws.Cells.Item($Row, 10).Value = $J
    $ws.Cells.Item($Row, 11).Value = $K
    $ws.Cells.Item($Row, 12).Value = $L
    $ws.Cells.Item($Row, 13).Value = $M
    $ws.Cells.Item($Row, 14).Value = $N
    $ws.Cells.Item($Row, 15).Value = $O
    $ws.Cells.Item($Row, 16).Value = $P
    $ws.Cells.Item($Row, 17).Value = $Q
    $ws.Cells.Item($Row, 18).Value = $R
    $ws.Cells.Item($Row, 19).Value = $S
    $ws.Cells.Item($Row, 20).Value = $T
}

# New row 699 - Packham's Triumph / Especial
Set-Row 699 8 "Terminal La Palmera de La Serena" "Coquimbo" 44747 4 "Fruta" `
    100104 "Frutos de pepita" 100104005 "Pera" "Packham's Triumph" "Especial" `
    20 210000 220000 215000 "`$/bins (450 kilos)" "Región de O'Higgins" 478 450

# New row 700 - Packham's Triumph / Primera
Set-Row 700 8 "Terminal La Palmera de La Serena" "Coquimbo" 44747 4 "Fruta" `
    100104 "Frutos de pepita" 100104005 "Pera" "Packham's Triumph" "Primera" `
    16 190000 200000 195000 "`$/bins (450 kilos)" "Región de O'Higgins" 433 450

# New row 701 - Packham's Triumph / Segunda
Set-Row 701 8 "Terminal La Palmera de La Serena" "Coquimbo" 44747 4 "Fruta" `
    100104 "Frutos de pepita" 100104005 "Pera" "Packham's Triumph" "Segunda" `
    14 160000 170000 165000 "`$/bins (450 kilos)" "Región de O'Higgins" 367 450
